# Auto-generated PowerShell Excel COM-interop script
# Applies the cryptos.xlsx price/volume/name/link updates described in the commit diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.605.42'
$ws.Range('E2').Value = '  +0.10%  '
$ws.Range('D3').Value = '1.847.25'
$ws.Range('E3').Value = '  -0.37%  '
$ws.Range('E4').Value = '  +0.33%  '
$ws.Range('D5').Value = '262.17'
$ws.Range('E5').Value = '  -1.34%  '
$ws.Range('E6').Value = '  +0.28%  '
$ws.Range('D7').Value = '0.5323'
$ws.Range('E7').Value = '  +1.75%  '
$ws.Range('D8').Value = '0.3166'
$ws.Range('E8').Value = '  -3.83%  '
$ws.Range('D9').Value = '0.06968'
$ws.Range('E9').Value = '  +2.14%  '
$ws.Range('D10').Value = '18.88'
$ws.Range('E10').Value = '  -0.09%  '
$ws.Range('D11').Value = '0.7729'
$ws.Range('E11').Value = '  -0.97%  '
$ws.Range('D12').Value = '0.07830'
$ws.Range('E12').Value = '  +0.53%  '
$ws.Range('D13').Value = '1.851.07'
$ws.Range('E13').Value = '  -0.08%  '
$ws.Range('D14').Value = '89.45'
$ws.Range('E14').Value = '  +0.95%  '
$ws.Range('D15').Value = '5.043'
$ws.Range('E15').Value = '  +0.20%  '
$ws.Range('B16').Value = 'BinanceUSD'
$ws.Range('C16').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D16').Value = '1.003'
$ws.Range('E16').Value = '  +0.41%  '
$ws.Range('B17').Value = 'Avalanche'
$ws.Range('C17').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D17').Value = '14.13'
$ws.Range('E17').Value = '  +0.90%  '
$ws.Range('D18').Value = '0.000007973'
$ws.Range('E18').Value = '  -0.25%  '
$ws.Range('D19').Value = '1.002'
$ws.Range('E19').Value = '  +0.25%  '
$ws.Range('D20').Value = '26.635.80'
$ws.Range('E20').Value = '  +0.18%  '
$ws.Range('B21').Value = 'Uniswap'
$ws.Range('C21').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D21').Value = '4.650'
$ws.Range('E21').Value = '  -0.06%  '
$ws.Range('B22').Value = 'Chainlink'
$ws.Range('C22').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D22').Value = '6.025'
$ws.Range('E22').Value = '  +0.48%  '
$ws.Range('B23').Value = 'Cosmos'
$ws.Range('C23').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D23').Value = '9.360'
$ws.Range('E23').Value = '  -2.28%  '
$ws.Range('B24').Value = 'Monero'
$ws.Range('C24').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D24').Value = '142.68'
$ws.Range('E24').Value = '  -1.53%  '
$ws.Range('B25').Value = 'LidoDAOToken'
$ws.Range('C25').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D25').Value = '2.216'
$ws.Range('E25').Value = '  -0.67%  '
$ws.Range('B26').Value = 'Toncoin'
$ws.Range('C26').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D26').Value = '1.695'
$ws.Range('E26').Value = '  +1.90%  '
$ws.Range('B27').Value = 'EthereumClassic'
$ws.Range('C27').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D27').Value = '17.13'
$ws.Range('E27').Value = '  +0.66%  '
$ws.Range('B28').Value = 'BitcoinCash'
$ws.Range('C28').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D28').Value = '111.68'
$ws.Range('E28').Value = '  -0.34%  '
$ws.Range('B29').Value = 'InternetComputer(DFINITY)'
$ws.Range('C29').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D29').Value = '4.325'
$ws.Range('E29').Value = '  +2.98%  '
$ws.Range('B30').Value = 'Stellar'
$ws.Range('C30').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D30').Value = '0.08778'
$ws.Range('E30').Value = '  +0.13%  '
$ws.Range('B31').Value = 'Filecoin'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D31').Value = '4.111'
$ws.Range('E31').Value = '  -1.07%  '
$ws.Range('B32').Value = 'Hedera'
$ws.Range('C32').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D32').Value = '0.04864'
$ws.Range('E32').Value = '  +0.18%  '
$ws.Range('B33').Value = 'ImmutableX'
$ws.Range('C33').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D33').Value = '0.7393'
$ws.Range('E33').Value = '  +2.37%  '
$ws.Range('B34').Value = 'ARBITRUM'
$ws.Range('C34').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D34').Value = '1.139'
$ws.Range('E34').Value = '  -0.28%  '
$ws.Range('B35').Value = 'HuobiToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D35').Value = '2.891'
$ws.Range('E35').Value = '  +1.34%  '
$ws.Range('B36').Value = 'MXToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D36').Value = '3.108'
$ws.Range('E36').Value = '  -0.02%  '
$ws.Range('B37').Value = 'RenderToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D37').Value = '2.370'
$ws.Range('E37').Value = '  +6.31%  '
$ws.Range('B38').Value = 'VeChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D38').Value = '0.01740'
$ws.Range('E38').Value = '  -2.70%  '
$ws.Range('B39').Value = 'TheSandbox'
$ws.Range('C39').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D39').Value = '0.4824'
$ws.Range('E39').Value = '  -1.63%  '
$ws.Range('B40').Value = 'TrustWalletToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D40').Value = '0.9068'
$ws.Range('E40').Value = '  -0.79%  '
$ws.Range('B41').Value = 'Quant'
$ws.Range('C41').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D41').Value = '109.18'
$ws.Range('E41').Value = '  -2.48%  '
$ws.Range('B42').Value = 'FraxShare'
$ws.Range('C42').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D42').Value = '5.918'
$ws.Range('E42').Value = '  -2.79%  '
$ws.Range('B43').Value = 'PaxDollar'
$ws.Range('C43').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D43').Value = '1.002'
$ws.Range('E43').Value = '  +0.28%  '
$ws.Range('B44').Value = 'Aptos'
$ws.Range('C44').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D44').Value = '7.706'
$ws.Range('E44').Value = '  -0.76%  '
$ws.Range('B45').Value = 'Decentraland'
$ws.Range('C45').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D45').Value = '0.4215'
$ws.Range('E45').Value = '  +0.31%  '
$ws.Range('B46').Value = 'EnergySwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D46').Value = '9.123'
$ws.Range('E46').Value = '  -0.34%  '
$ws.Range('B47').Value = 'Algorand'
$ws.Range('C47').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D47').Value = '0.1251'
$ws.Range('E47').Value = '  +0.39%  '
$ws.Range('B48').Value = 'Elrond'
$ws.Range('C48').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D48').Value = '35.10'
$ws.Range('E48').Value = '  -0.35%  '
$ws.Range('B49').Value = 'Cronos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D49').Value = '0.05814'
$ws.Range('E49').Value = '  -2.30%  '
$ws.Range('B50').Value = 'EOS'
$ws.Range('C50').Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range('D50').Value = '0.8987'
$ws.Range('E50').Value = '  +0.76%  '
$ws.Range('B51').Value = 'Aave'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D51').Value = '60.45'
$ws.Range('E51').Value = '  +0.52%  '
